$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "31.291.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.04%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.989.06"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +6.28%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9981"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.8163"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +73.40%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "253.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9971"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.30%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3462"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +20.15%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.71"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +17.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07027"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.75%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8427"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +16.17%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08107"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.25%  "

$ws.Range("B13").Value = "Litecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "100.89"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.87%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.987.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.24%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.526"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.62%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "273.58"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.32%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "31.292.96"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.08%  "

$ws.Range("E18").Value = "  +7.20%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007938"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.89%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.800"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +10.46%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.248.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.48%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9983"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.17%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9963"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.38%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.967"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +11.74%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.779"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1509"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +56.64%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "164.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.28%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.79%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.202"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +17.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.569"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.78%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.567"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.87%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.345"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.86%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.326"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.10%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05183"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.63%  "

$ws.Range("E35").Value = "  +8.27%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7590"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.96%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.756"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.66%  "

$ws.Range("E38").Value = "  +5.69%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.906"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.43%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.607"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.85%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4712"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +11.47%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "78.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.22%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.090"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.28%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8536"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.64%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.85%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9969"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.24%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.936"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.29%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.510"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.92%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4321"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.17%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.75"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.04%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1196"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +12.67%  "
